$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for specific rows per repulled data / mean calculation
$ws.Range("F2").Value = -2
$ws.Range("F7").Value = -1
$ws.Range("F8").Value = 1
$ws.Range("F9").Value = -6
$ws.Range("F10").Value = -2
$ws.Range("F11").Value = -6
$ws.Range("F13").Value = -3
$ws.Range("F14").Value = -1
$ws.Range("F17").Value = 6
$ws.Range("F19").Value = 2
$ws.Range("F21").Value = -2
$ws.Range("F23").Value = -10
$ws.Range("F24").Value = -2
